$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet "Provider Submissions" -> "Data Quality"
$ws.Name = "Data Quality"

# Turn on AutoFilter over the header row B12:K12
$ws.Range("B12:K12").AutoFilter()

# Register the hidden _FilterDatabase defined name scoped to this sheet,
# as Excel does automatically when AutoFilter is applied via the UI.
$filterDbName = $ws.Names.Add("_xlnm._FilterDatabase", "='Data Quality'!`$B`$12:`$K`$12")
$filterDbName.Visible = $false

# Update the current selection to the filtered header range
$ws.Range("B12:K12").Select()
